$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text with the new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 2.87 = 10646.55 pesos`n✅ 10646.55 pesos = 2.86 = 962.11 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 348
$ws2.Range("O10").Value = 3705
$ws2.Range("N12").Value = 3717.03
$ws2.Range("O12").Value = 335.9
